# Update the fixed "datetimeFigureOut" date field text from 5/10/2011 to
# 5/11/2011 everywhere it appears: the slide master, every slide layout
# (custom layout), and the notes master.

$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "5/10/2011") {
                $tr.Text = "5/11/2011"
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every layout belonging to the slide master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DateShape $notesMaster.Shapes
